# Update "想去人数" (want-to-go count) and "最低票价" (minimum ticket price)
# figures across the three data sheets: 展览 (Exhibition), 演出 (Performance)
# and 全部类型 (All Types, which mirrors rows from the first two sheets).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 621
$ws1.Range("F4").Value = 585
$ws1.Range("G4").Value = 50
$ws1.Range("F5").Value = 531
$ws1.Range("F6").Value = 299
$ws1.Range("F7").Value = 2693
$ws1.Range("F8").Value = 459
$ws1.Range("F9").Value = 7479
$ws1.Range("F10").Value = 197
$ws1.Range("F11").Value = 460
$ws1.Range("F12").Value = 27
$ws1.Range("F13").Value = 243
$ws1.Range("F14").Value = 41

# --- 演出 (Performance) sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 18

# --- 全部类型 (All Types) sheet, which repeats the same rows ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 621
$ws4.Range("F4").Value = 585
$ws4.Range("G4").Value = 50
$ws4.Range("F5").Value = 531
$ws4.Range("F6").Value = 299
$ws4.Range("F8").Value = 18
$ws4.Range("F9").Value = 2693
$ws4.Range("F10").Value = 459
$ws4.Range("F11").Value = 7479
$ws4.Range("F12").Value = 197
$ws4.Range("F13").Value = 460
$ws4.Range("F14").Value = 27
$ws4.Range("F17").Value = 243
$ws4.Range("F18").Value = 41
